# Updates cryptos list price/volume data (and reorders Mantle/WhiteBITCoin rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-44: update Price (D) and Volume(1h) (E) text values ---
$ws.Range("D2").Value = "'59.305.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.19%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.522.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.43%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'536.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.05%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'140.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.40%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.26%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.563"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.95%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.529.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.48%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0993"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.61%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.73%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'5.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.96%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.44%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.967.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.40%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'23.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.34%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'59.216.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.36%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.0000141"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.80%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.537.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.41%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'10.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.93%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'4.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.54%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'320.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.25%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.14%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.00%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'61.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.78%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.420"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.96%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.166"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.40%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.997"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.49%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.15%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -2.29%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0769"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.31%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.44%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'162.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.91%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.29%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -9.80%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.53%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'18.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.52%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'4.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -5.11%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.52%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'36.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.39%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'3.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.32%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'5.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -9.90%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'288.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -6.58%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.805"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.85%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.24%  "
$ws.Range("E44").Style = "Normal"

# --- Rows 45-46: Mantle and WhiteBITCoin swap places, with updated values ---
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.599"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.50%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "'10.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.82%  "
$ws.Range("E46").Style = "Normal"

# --- Rows 47-51: update Price (D) and Volume(1h) (E) text values ---
$ws.Range("D47").Value = "'124.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.02%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0926"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.63%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'18.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.30%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0508"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.86%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0224"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.96%  "
$ws.Range("E51").Style = "Normal"
